# "update data test for web, mobile, and api"
# Replace the "Archie / arch* " sample data in Sheet1 with the new "Joong" data set.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- fullname column (C) ---
$ws.Range("C3:C10").Value = "Joong A"

# --- email column (D) ---
$ws.Range("D3").Value = "joongarc"
$ws.Range("D4:D9").Value = "joongarc@gmail.com"
$ws.Range("D10").Value = "joong21@gmail.com"

# --- password column (E) ---
$ws.Range("E3").Value = "Joong123!"
$ws.Range("E4").Value = "Jng12!"
$ws.Range("E5").Value = "Joongjoong12345!"
$ws.Range("E6").Value = "joong123!"
$ws.Range("E7").Value = "Joong123"
$ws.Range("E8").Value = "Joongs#!"
$ws.Range("E9").Value = "Joong123!"
# E10 keeps its existing value "Joong123!"

# --- hyperlinks: re-point the mailto links at the new e-mail addresses ---
$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:joong21@gmail.com") | Out-Null
$ws.Range("D10").Font.Underline = $true

$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:joongarc@gmail.com") | Out-Null
$ws.Range("D4").Font.Underline = $true

$ws.Hyperlinks.Add($ws.Range("D5:D9"), "mailto:joongarc@gmail.com", "", "", "joongarc@gmail.com") | Out-Null
$ws.Range("D5:D9").Font.Underline = $true

# --- column widths: fullname / email columns got wider ---
$ws.Columns.Item(3).ColumnWidth = 14.43
$ws.Columns.Item(4).ColumnWidth = 20.43

$wb.Save()
